# Apply odds updates to the single worksheet ("Sheet1") of the workbook.
# Each entry maps a cell reference to its new numeric value, matching the
# target diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "N2"   = 9
    "H11"  = 3.6
    "J11"  = 5.1
    "L11"  = 2.2
    "O11"  = 1.32
    "P11"  = 3.1
    "W11"  = 11.75
    "X11"  = 27
    "Y11"  = 16
    "Z11"  = 90
    "AE11" = 17
    "AG11" = 800
    "AH11" = 6.4
    "AI11" = 7.6
    "AL11" = 13.5
    "AM11" = 28
    "AN11" = 6.5
    "AO11" = 28
    "AP11" = 35
    "AQ11" = 175
    "AS11" = 450
    "AU11" = 7.6
    "AW11" = 3.5
    "AX11" = 8
    "BA11" = 55
    "M13"  = 1.1
    "N13"  = 7
    "Q13"  = 2.4
    "R13"  = 1.53
    "G15"  = 2.25
    "I15"  = 3.2
    "L15"  = 3.75
    "N15"  = 9.5
    "Q15"  = 2.03
    "R15"  = 1.83
    "AC15" = 9.5
    "AI15" = 17
    "BB15" = 201
    "Q16"  = 2.15
    "R16"  = 1.67
}

foreach ($cellRef in $changes.Keys) {
    $ws.Range($cellRef).Value = $changes[$cellRef]
}
